$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.742.63"
$ws.Range("E2").Value = "  +6.83%  "
$ws.Range("D3").Value = "2.312.12"
$ws.Range("E3").Value = "  +5.89%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.18%  "
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  +11.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.35%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "2.663.85"
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("D15").Value = "2.307.92"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.819"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.33%  "
$ws.Range("D18").Value = "46.693.18"
$ws.Range("E18").Value = "  +7.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +25.05%  "
$ws.Range("E20").Value = "  +8.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.96%  "
$ws.Range("E25").Value = "  +10.68%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "43.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +23.90%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.56%  "
$ws.Range("E30").Value = "  +7.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +11.10%  "
$ws.Range("E34").Value = "  +6.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.70%  "
$ws.Range("E36").Value = "  +13.32%  "
$ws.Range("E37").Value = "  +3.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +20.64%  "
$ws.Range("E40").Value = "  +16.72%  "
$ws.Range("E41").Value = "  +14.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0307"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.50%  "
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +22.68%  "
$ws.Range("D45").Value = "1.839.91"
$ws.Range("E45").Value = "  +6.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +25.03%  "
$ws.Range("E47").Value = "  +18.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.79%  "
$ws.Range("E49").Value = "  +13.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.95%  "
